$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-FirstInlineShape($range, $newName) {
    $ish = $range.InlineShapes
    if ($ish.Count -gt 0) {
        $shape = $ish.Item(1)
        # Re-scope to the shape's own (narrow) range before renaming so the
        # COM host re-anchors on the run that actually holds the drawing --
        # this matters when the picture isn't in the first paragraph of the
        # header/footer story.
        $narrow = $shape.Range.InlineShapes
        $narrow.Item(1).Name = $newName
    }
}

# Headers -- BTec_Logo-Orange picture: image2.jpg -> image1.jpg
for ($i = 1; $i -le 2; $i++) {
    $h = $sec.Headers.Item($i)
    if ($h.Exists) {
        Rename-FirstInlineShape $h.Range "image1.jpg"
    }
}

# Footers -- Pearson logo picture: image1.png -> image2.png
for ($i = 1; $i -le 2; $i++) {
    $f = $sec.Footers.Item($i)
    if ($f.Exists) {
        Rename-FirstInlineShape $f.Range "image2.png"
    }
}
